# Update the arithmetic problems ("XX÷Y=") in the single 20x5 practice
# table. Replacements are targeted per table cell (row, column) so that
# duplicate source strings (e.g. "30÷2=" appears twice) are each mapped
# to their own distinct replacement.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$changes = @(
    @{ Row = 1;  Col = 1; Old = "16÷6="; New = "68÷6=" },
    @{ Row = 1;  Col = 2; Old = "11÷5="; New = "68÷2=" },
    @{ Row = 1;  Col = 3; Old = "94÷2="; New = "81÷8=" },
    @{ Row = 1;  Col = 4; Old = "30÷2="; New = "76÷5=" },
    @{ Row = 1;  Col = 5; Old = "57÷7="; New = "67÷6=" },

    @{ Row = 5;  Col = 1; Old = "12÷6="; New = "21÷2=" },
    @{ Row = 5;  Col = 2; Old = "79÷2="; New = "61÷5=" },
    @{ Row = 5;  Col = 3; Old = "93÷2="; New = "88÷4=" },
    @{ Row = 5;  Col = 4; Old = "53÷8="; New = "62÷8=" },
    @{ Row = 5;  Col = 5; Old = "30÷2="; New = "80÷5=" },

    @{ Row = 9;  Col = 1; Old = "54÷5="; New = "74÷7=" },
    @{ Row = 9;  Col = 2; Old = "59÷9="; New = "88÷4=" },
    @{ Row = 9;  Col = 3; Old = "57÷9="; New = "64÷3=" },
    @{ Row = 9;  Col = 4; Old = "36÷8="; New = "97÷8=" },
    @{ Row = 9;  Col = 5; Old = "77÷9="; New = "84÷9=" },

    @{ Row = 13; Col = 1; Old = "75÷2="; New = "54÷3=" },
    @{ Row = 13; Col = 2; Old = "48÷5="; New = "16÷8=" },
    @{ Row = 13; Col = 3; Old = "81÷7="; New = "68÷5=" },
    @{ Row = 13; Col = 4; Old = "83÷6="; New = "92÷3=" },
    @{ Row = 13; Col = 5; Old = "93÷9="; New = "79÷3=" },

    @{ Row = 17; Col = 1; Old = "49÷9="; New = "12÷8=" },
    @{ Row = 17; Col = 2; Old = "83÷4="; New = "28÷4=" },
    @{ Row = 17; Col = 3; Old = "30÷4="; New = "79÷3=" },
    @{ Row = 17; Col = 4; Old = "46÷7="; New = "90÷4=" },
    @{ Row = 17; Col = 5; Old = "65÷7="; New = "39÷5=" }
)

foreach ($change in $changes) {
    $cellRange = $t.Cell($change.Row, $change.Col).Range
    $cellRange.Text = $change.New
}
